$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column C - this shifts the existing
# "OS" (old C) and "State" (old D) columns one to the right, becoming D and E,
# and makes room for the new "Model Name" column at C.
$ws.Columns("C:C").Insert()
$ws.Columns("C:C").ColumnWidth = 16.42578125

# Row 6 (new device): Samsung SM-G920F (Galaxy S6), Android 7
$ws.Range("A6").Value = "Samsung"
$ws.Range("B6").Value = "SM-G920F"
$ws.Range("D6").Value = "Android 7"
$ws.Range("E6").Value = "WIFI MAC not possible"
$ws.Range("E6").WrapText = $true

# Row 7 (new device): Motorala XT-1925-5 (moto g6), Android 9
$ws.Range("A7").Value = "Motorala"
$ws.Range("D7").Value = "Android 9"
$ws.Range("B7").Value = "XT-1925-5"

# New "Model Name" header/column
$ws.Range("C4").Value = "Model Name"
$ws.Range("C6").Value = "Galaxy S6"
$ws.Range("C7").Value = "moto g6"
$ws.Range("C5").Value = "Galaxy Tab 3"

$ws.Range("E7").Value = "Everythings works"

$ws.Range("E8").Select()
